$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF (column I) values for rows 30 through 48 from the
# previous value (24.44214285714286) to the new value (18.473)
$ws.Range("I30:I48").Value = 18.473
